$wb = $excel.ActiveWorkbook

# Target the DatosCuenta sheet (first sheet), which is the active/selected sheet
$ws = $wb.Worksheets.Item("DatosCuenta")
$ws.Activate()

# Update the text value shared by A2 and B2 from "SmokeOcho" to "SmokeDiez"
$ws.Range("A2").Value = "SmokeDiez"
$ws.Range("B2").Value = "SmokeDiez"

# Update numeric values in C2 and D2
$ws.Range("C2").Value = 21546910
$ws.Range("D2").Value = 142

# Update the active selection to C2
$ws.Range("C2").Select()

$wb.Save()
